# Updated get balance button
# Populate the crypto balance table with the full list of currencies
# returned by the "get balance" call, replacing the old placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (currency GAL / balance 2.57941800) is left untouched.

# Extend the bordered/bold index-column formatting (column A, style of A2:A4)
# down through the new rows before writing their values.
$ws.Range("A4").Copy()
$ws.Range("A5:A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 3-11: new currencies with numeric balances returned by the balance call.
$data = @(
    @{ Row = 3;  Idx = 1; Name = "SHIB"; Balance = 1307702.99 },
    @{ Row = 4;  Idx = 2; Name = "FIRO"; Balance = 0.08190859 },
    @{ Row = 5;  Idx = 3; Name = "SOL";  Balance = 0.36963 },
    @{ Row = 6;  Idx = 4; Name = "BRL";  Balance = 1.14509252 },
    @{ Row = 7;  Idx = 5; Name = "ADA";  Balance = 33.7662 },
    @{ Row = 8;  Idx = 6; Name = "XRP";  Balance = 172.790793 },
    @{ Row = 9;  Idx = 7; Name = "USDT"; Balance = 0.0022175 },
    @{ Row = 10; Idx = 8; Name = "ETH";  Balance = 0.00945705 },
    @{ Row = 11; Idx = 9; Name = "BTC";  Balance = 0.0000061 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Idx
    $ws.Cells.Item($r, 2).Value = $item.Name
    $ws.Cells.Item($r, 3).Value = $item.Balance
}
